$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 50: finalize the STR buy — mark DONE, stamp a finalized date, add a fee
# ---------------------------------------------------------------------------
$ws.Range("H50").Value = "DONE"

$ws.Range("I50").Value = 42863.590324074074
$ws.Range("A50").Copy()
$ws.Range("I50").PasteSpecial(-4122)

$ws.Range("J50").Value = "1.20900000 STR (0.15%)"

# ---------------------------------------------------------------------------
# Row 51: a brand-new STR sell that shows up 10 minutes later
# ---------------------------------------------------------------------------
$ws.Range("A51").Value = 42863.686921296299
$ws.Range("A50").Copy()
$ws.Range("A51").PasteSpecial(-4122)

# B51 needs the rich-text "Sell" run (leading spaces + red "Sell") already
# used elsewhere in the sheet (e.g. B30) — copy it so the shared string is
# reused instead of duplicated.
$ws.Range("B30").Copy()
$ws.Range("B51").PasteSpecial(-4163)

$ws.Range("C50").Copy()
$ws.Range("C51").PasteSpecial(-4163)

# D51 is a numeric-looking text value ("          0.05218998" + newline).
# Force text formatting first so COM doesn't coerce it to a Double, then
# restamp the proper cell style (matches D50's wrap-text style) afterwards.
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "          0.05218998`n"
$ws.Range("D50").Copy()
$ws.Range("D51").PasteSpecial(-4122)

$ws.Range("E51").Value = "         0.063USDT"

$ws.Range("F50").Copy()
$ws.Range("F51").PasteSpecial(-4163)

$ws.Range("G50").Copy()
$ws.Range("G51").PasteSpecial(-4163)

$ws.Range("H51").Value = "IN PROGRESS"

$ws.Range("I51").Value = 42863.590324074074
$ws.Range("I50").Copy()
$ws.Range("I51").PasteSpecial(-4122)

$ws.Range("K50").Copy()
$ws.Range("K51").PasteSpecial(-4163)

# Match row height used throughout the log table.
$ws.Rows.Item(51).RowHeight = 14.25

# ---------------------------------------------------------------------------
# View state: scrolled down and selection moved to where data entry continued
# ---------------------------------------------------------------------------
$ws.Range("E54").Select() | Out-Null

$excel.CutCopyMode = $false
